$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.366.89'
$ws.Range("E2").Value = '  +1.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.899.08'
$ws.Range("E3").Value = '  +2.82%  '
$ws.Range("E4").Value = '  +0.50%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.46'
$ws.Range("E5").Value = '  +2.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.648'
$ws.Range("E6").Value = '  +4.68%  '
$ws.Range("E7").Value = '  +0.62%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.44'
$ws.Range("E8").Value = '  -2.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.342'
$ws.Range("E9").Value = '  +4.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '49.87'
$ws.Range("E10").Value = '  +7.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0709'
$ws.Range("E11").Value = '  +2.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0995'
$ws.Range("E12").Value = '  +0.77%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.175.87'
$ws.Range("E13").Value = '  +2.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '12.09'
$ws.Range("E14").Value = '  +6.56%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.692'
$ws.Range("E15").Value = '  +2.92%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.900.30'
$ws.Range("E16").Value = '  +2.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.83'
$ws.Range("E17").Value = '  +1.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '35.372.13'
$ws.Range("E18").Value = '  +1.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.54'
$ws.Range("E19").Value = '  +2.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0812'
$ws.Range("E20").Value = '  +2.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '241.30'
$ws.Range("E21").Value = '  +0.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.48'
$ws.Range("E22").Value = '  +3.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.74'
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.01'
$ws.Range("E24").Value = '  +0.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.30'
$ws.Range("E25").Value = '  +1.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.28'
$ws.Range("E26").Value = '  +23.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.05'
$ws.Range("E27").Value = '  +0.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.32'
$ws.Range("E28").Value = '  +4.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.16'
$ws.Range("E29").Value = '  +3.43%  '
$ws.Range("E30").Value = '  +1.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.11'
$ws.Range("E31").Value = '  +3.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0562'
$ws.Range("E32").Value = '  +1.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.02'
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.923'
$ws.Range("E34").Value = '  +18.06%  '
$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.10'
$ws.Range("E35").Value = '  +2.75%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.74'
$ws.Range("E36").Value = '  +4.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.04'
$ws.Range("E37").Value = '  +2.17%  '
$ws.Range("E38").Value = '  +2.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0209'
$ws.Range("E39").Value = '  +3.59%  '
$ws.Range("E40").Value = '  +1.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0631'
$ws.Range("E41").Value = '  +13.91%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '15.70'
$ws.Range("E42").Value = '  +5.79%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '89.15'
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.337.04'
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.35'
$ws.Range("E45").Value = '  +1.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '47.59'
$ws.Range("E46").Value = '  +39.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.41'
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("E48").Value = '  +1.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '12.36'
$ws.Range("E49").Value = '  -14.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.49'
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.082.26'
$ws.Range("E51").Value = '  +2.67%  '
